# First run of pipeline.
# Relabel the "Unprocessed <source>" / "Resources - Wind" rows (rows 2-11,
# column D and C11) to their new "[from Resources]" / "[of Wind]" wording,
# and move the active selection to C11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "Coal & Coal products [from Resources]"
$ws.Range("D3").Value  = "Oil & oil products [from Resources]"
$ws.Range("D4").Value  = "Natural gas [from Resources]"
$ws.Range("D5").Value  = "Non-specified primary biofuels and waste [from Resources]"
$ws.Range("D6").Value  = "Nuclear [from Resources]"
$ws.Range("D7").Value  = "Hydro [from Resources]"
$ws.Range("D8").Value  = "Geothermal [from Resources]"
$ws.Range("D9").Value  = "Solar photovoltaics [from Resources]"
$ws.Range("D10").Value = "Solar thermal [from Resources]"
$ws.Range("C11").Value = "Resources [of Wind]"
$ws.Range("D11").Value = "Wind [from Resources]"

$ws.Range("C11").Select()
